$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.6591047123919509
$ws.Cells.Item(2, 3).Value = 0.06979013589223371
$ws.Cells.Item(2, 4).Value = 0.0764930503263459
$ws.Cells.Item(2, 5).Value = 0.1128240276632582
$ws.Cells.Item(2, 6).Value = 1.807478639038322
$ws.Cells.Item(2, 8).Value = 0.07973214163530429
$ws.Cells.Item(2, 11).Value = 0.651558464049657
$ws.Cells.Item(2, 13).Value = 0.3084085407963499
$ws.Cells.Item(2, 14).Value = 2.635021442869999
$ws.Cells.Item(3, 2).Value = 0.6141290824258192
$ws.Cells.Item(3, 3).Value = 0.06204465075281007
$ws.Cells.Item(3, 4).Value = 0.07671788444945538
$ws.Cells.Item(3, 5).Value = 0.1034124664791491
$ws.Cells.Item(3, 6).Value = 1.775948807909614
$ws.Cells.Item(3, 8).Value = 0.07973214163530429
$ws.Cells.Item(3, 11).Value = 0.601413058193998
$ws.Cells.Item(3, 13).Value = 0.2838194162092975
$ws.Cells.Item(3, 14).Value = 2.632067726334895
$ws.Cells.Item(4, 2).Value = 0.5869108781935779
$ws.Cells.Item(4, 3).Value = 0.05731793280595809
$ws.Cells.Item(4, 4).Value = 0.07685709186201173
$ws.Cells.Item(4, 5).Value = 0.09770089020240391
$ws.Cells.Item(4, 6).Value = 1.757449805792177
$ws.Cells.Item(4, 8).Value = 0.07973214163530429
$ws.Cells.Item(4, 11).Value = 0.5710057226790184
$ws.Cells.Item(4, 13).Value = 0.2689034936060395
$ws.Cells.Item(4, 14).Value = 2.630894410618922
$ws.Cells.Item(5, 2).Value = 0.5759186639842824
$ws.Cells.Item(5, 3).Value = 0.05539892876647912
$ws.Cells.Item(5, 4).Value = 0.07691410927512621
$ws.Cells.Item(5, 5).Value = 0.09538999018739247
$ws.Cells.Item(5, 6).Value = 1.750127065168229
$ws.Cells.Item(5, 8).Value = 0.07973214163530429
$ws.Cells.Item(5, 11).Value = 0.5587100363281934
$ws.Cells.Item(5, 13).Value = 0.2628704109428028
$ws.Cells.Item(5, 14).Value = 2.630576868961612
$ws.Cells.Item(6, 2).Value = 0.5740994085061004
$ws.Cells.Item(6, 3).Value = 0.05508070898170558
$ws.Cells.Item(6, 4).Value = 0.07692359438741647
$ws.Cells.Item(6, 5).Value = 0.09500726244970537
$ws.Cells.Item(6, 6).Value = 1.748924142906333
$ws.Cells.Item(6, 8).Value = 0.07973214163530429
$ws.Cells.Item(6, 11).Value = 0.5566741022619794
$ws.Cells.Item(6, 13).Value = 0.2618713433420226
$ws.Cells.Item(6, 14).Value = 2.630533829007646
$ws.Cells.Item(7, 2).Value = 0.586762231313827
$ws.Cells.Item(7, 3).Value = 0.05729202357916563
$ws.Cells.Item(7, 4).Value = 0.07685785965097747
$ws.Cells.Item(7, 5).Value = 0.09766965768806557
$ws.Cells.Item(7, 6).Value = 1.757350175994475
$ws.Cells.Item(7, 8).Value = 0.07973214163530429
$ws.Cells.Item(7, 11).Value = 0.5708395125882078
$ws.Cells.Item(7, 13).Value = 0.2688219466259838
$ws.Cells.Item(7, 14).Value = 2.630889478464752
$ws.Cells.Item(8, 2).Value = 0.6435144523371434
$ws.Cells.Item(8, 3).Value = 0.0671133837222726
$ws.Cells.Item(8, 4).Value = 0.0765703315062396
$ws.Cells.Item(8, 5).Value = 0.1095648137323337
$ws.Cells.Item(8, 6).Value = 1.796428184714287
$ws.Cells.Item(8, 8).Value = 0.07973214163530429
$ws.Cells.Item(8, 11).Value = 0.6341886674001955
$ws.Cells.Item(8, 13).Value = 0.2998921696832753
$ws.Cells.Item(8, 14).Value = 2.633869860455675
$ws.Cells.Item(9, 2).Value = 0.757980601401357
$ws.Cells.Item(9, 3).Value = 0.08661038358275164
$ws.Cells.Item(9, 4).Value = 0.07601575579693787
$ws.Cells.Item(9, 5).Value = 0.1334376943439537
$ws.Cells.Item(9, 6).Value = 1.879920921771259
$ws.Cells.Item(9, 8).Value = 0.07973214163530429
$ws.Cells.Item(9, 11).Value = 0.7614801266468305
$ws.Cells.Item(9, 13).Value = 0.3622883877112244
$ws.Cells.Item(9, 14).Value = 2.644815051827663
$ws.Cells.Item(10, 2).Value = 0.8440589803889225
$ws.Cells.Item(10, 3).Value = 0.1010905547000505
$ws.Cells.Item(10, 4).Value = 0.07561403926781285
$ws.Cells.Item(10, 5).Value = 0.1513309495337865
$ws.Cells.Item(10, 6).Value = 1.945500955783601
$ws.Cells.Item(10, 8).Value = 0.07973214163530429
$ws.Cells.Item(10, 11).Value = 0.8569235279089469
$ws.Cells.Item(10, 13).Value = 0.4090641847685674
$ws.Cells.Item(10, 14).Value = 2.655996909133023
$ws.Cells.Item(11, 2).Value = 0.883658487502089
$ws.Cells.Item(11, 3).Value = 0.1077142584296382
$ws.Cells.Item(11, 4).Value = 0.07543255222692746
$ws.Cells.Item(11, 5).Value = 0.1595524870685381
$ws.Cells.Item(11, 6).Value = 1.976267830655786
$ws.Cells.Item(11, 8).Value = 0.07973214163530429
$ws.Cells.Item(11, 11).Value = 0.9007731225269424
$ws.Cells.Item(11, 13).Value = 0.4305550452683065
$ws.Cells.Item(11, 14).Value = 2.661772659867268
$ws.Cells.Item(12, 2).Value = 0.8987179914255989
$ws.Cells.Item(12, 3).Value = 0.110227917974612
$ws.Cells.Item(12, 4).Value = 0.07536401123984859
$ws.Cells.Item(12, 5).Value = 0.162677871343945
$ws.Cells.Item(12, 6).Value = 1.988053640294908
$ws.Cells.Item(12, 8).Value = 0.07973214163530429
$ws.Cells.Item(12, 11).Value = 0.9174407350377862
$ws.Cells.Item(12, 13).Value = 0.4387242121703849
$ws.Cells.Item(12, 14).Value = 2.664059396719466
$ws.Cells.Item(13, 2).Value = 0.8954718025055968
$ws.Cells.Item(13, 3).Value = 0.1096863136490356
$ws.Cells.Item(13, 4).Value = 0.07537876451428627
$ws.Cells.Item(13, 5).Value = 0.1620042224060256
$ws.Cells.Item(13, 6).Value = 1.98550933671379
$ws.Cells.Item(13, 8).Value = 0.07973214163530429
$ws.Cells.Item(13, 11).Value = 0.9138482669619066
$ws.Cells.Item(13, 13).Value = 0.4369634469569235
$ws.Cells.Item(13, 14).Value = 2.663562469928877
$ws.Cells.Item(14, 2).Value = 0.8848961559837676
$ws.Cells.Item(14, 3).Value = 0.1079209495268287
$ws.Cells.Item(14, 4).Value = 0.07542690960459097
$ws.Cells.Item(14, 5).Value = 0.1598093707382304
$ws.Cells.Item(14, 6).Value = 1.977234745268134
$ws.Cells.Item(14, 8).Value = 0.07973214163530429
$ws.Cells.Item(14, 11).Value = 0.9021431164280784
$ws.Cells.Item(14, 13).Value = 0.4312265024039164
$ws.Cells.Item(14, 14).Value = 2.66195879236983
$ws.Cells.Item(15, 2).Value = 0.8784266197642978
$ws.Cells.Item(15, 3).Value = 0.1068403214291607
$ws.Cells.Item(15, 4).Value = 0.07545642396582863
$ws.Cells.Item(15, 5).Value = 0.1584665411666819
$ws.Cells.Item(15, 6).Value = 1.972183929869118
$ws.Cells.Item(15, 8).Value = 0.07973214163530429
$ws.Cells.Item(15, 11).Value = 0.8949815604991045
$ws.Cells.Item(15, 13).Value = 0.4277165164709658
$ws.Cells.Item(15, 14).Value = 2.660989477408776
$ws.Cells.Item(16, 2).Value = 0.8414799938523174
$ws.Cells.Item(16, 3).Value = 0.1006584288193437
$ws.Cells.Item(16, 4).Value = 0.07562592544726066
$ws.Cells.Item(16, 5).Value = 0.1507953246818232
$ws.Cells.Item(16, 6).Value = 1.943509128293954
$ws.Cells.Item(16, 8).Value = 0.07973214163530429
$ws.Cells.Item(16, 11).Value = 0.8540665972286661
$ws.Cells.Item(16, 13).Value = 0.4076640225996258
$ws.Cells.Item(16, 14).Value = 2.655633358498278
$ws.Cells.Item(17, 2).Value = 0.8189279777695333
$ws.Cells.Item(17, 3).Value = 0.0968755139171833
$ws.Cells.Item(17, 4).Value = 0.07573023327228512
$ws.Cells.Item(17, 5).Value = 0.1461104487950777
$ws.Cells.Item(17, 6).Value = 1.926157848115054
$ws.Cells.Item(17, 8).Value = 0.07973214163530429
$ws.Cells.Item(17, 11).Value = 0.8290776494207535
$ws.Cells.Item(17, 13).Value = 0.3954172098195912
$ws.Cells.Item(17, 14).Value = 2.652524392397012
$ws.Cells.Item(18, 2).Value = 0.8059982311525573
$ws.Cells.Item(18, 3).Value = 0.09470311771480056
$ws.Cells.Item(18, 4).Value = 0.07579034634720916
$ws.Cells.Item(18, 5).Value = 0.1434235137650717
$ws.Cells.Item(18, 6).Value = 1.916265736855237
$ws.Cells.Item(18, 8).Value = 0.07973214163530429
$ws.Cells.Item(18, 11).Value = 0.814745283887504
$ws.Cells.Item(18, 13).Value = 0.3883931201024851
$ws.Cells.Item(18, 14).Value = 2.650801025497543
$ws.Cells.Item(19, 2).Value = 0.8016275678141938
$ws.Cells.Item(19, 3).Value = 0.09396816755042892
$ws.Cells.Item(19, 4).Value = 0.07581071975580933
$ws.Cells.Item(19, 5).Value = 0.142515073891083
$ws.Cells.Item(19, 6).Value = 1.9129315105057
$ws.Cells.Item(19, 8).Value = 0.07973214163530429
$ws.Cells.Item(19, 11).Value = 0.8098995485277669
$ws.Cells.Item(19, 13).Value = 0.3860182954256146
$ws.Cells.Item(19, 14).Value = 2.650228642281206
$ws.Cells.Item(20, 2).Value = 0.8213243720045682
$ws.Cells.Item(20, 3).Value = 0.09727785497670993
$ws.Cells.Item(20, 4).Value = 0.07571911729243652
$ws.Cells.Item(20, 5).Value = 0.1466083646195315
$ws.Cells.Item(20, 6).Value = 1.927995820423348
$ws.Cells.Item(20, 8).Value = 0.07973214163530429
$ws.Cells.Item(20, 11).Value = 0.8317335575214315
$ws.Cells.Item(20, 13).Value = 0.3967188351723152
$ws.Cells.Item(20, 14).Value = 2.652848634129455
$ws.Cells.Item(21, 2).Value = 0.88800073819408
$ws.Cells.Item(21, 3).Value = 0.1084393318733419
$ws.Cells.Item(21, 4).Value = 0.0754127631954491
$ws.Cells.Item(21, 5).Value = 0.1604537219680822
$ws.Cells.Item(21, 6).Value = 1.979661521970286
$ws.Cells.Item(21, 8).Value = 0.07973214163530429
$ws.Cells.Item(21, 11).Value = 0.9055794956546777
$ws.Cells.Item(21, 13).Value = 0.4329107354226238
$ws.Cells.Item(21, 14).Value = 2.662427124620194
$ws.Cells.Item(22, 2).Value = 0.9319510374986635
$ws.Cells.Item(22, 3).Value = 0.1157655844970975
$ws.Cells.Item(22, 4).Value = 0.07521361868668386
$ws.Cells.Item(22, 5).Value = 0.1695729386689493
$ws.Cells.Item(22, 6).Value = 2.01421563327375
$ws.Cells.Item(22, 8).Value = 0.07973214163530429
$ws.Cells.Item(22, 11).Value = 0.9542079734535207
$ws.Cells.Item(22, 13).Value = 0.4567454280978751
$ws.Cells.Item(22, 14).Value = 2.669267917465845
$ws.Cells.Item(23, 2).Value = 0.9084596088641206
$ws.Cells.Item(23, 3).Value = 0.1118524879337031
$ws.Cells.Item(23, 4).Value = 0.07531980618042766
$ws.Cells.Item(23, 5).Value = 0.1646992916983265
$ws.Cells.Item(23, 6).Value = 1.995701146475142
$ws.Cells.Item(23, 8).Value = 0.07973214163530429
$ws.Cells.Item(23, 11).Value = 0.9282203552013186
$ws.Cells.Item(23, 13).Value = 0.4440076534438475
$ws.Cells.Item(23, 14).Value = 2.665563557624438
$ws.Cells.Item(24, 2).Value = 0.8202408509839074
$ws.Cells.Item(24, 3).Value = 0.09709594895801388
$ws.Cells.Item(24, 4).Value = 0.07572414237469083
$ws.Cells.Item(24, 5).Value = 0.1463832367711291
$ws.Cells.Item(24, 6).Value = 1.927164613501532
$ws.Cells.Item(24, 8).Value = 0.07973214163530429
$ws.Cells.Item(24, 11).Value = 0.8305327152614836
$ws.Cells.Item(24, 13).Value = 0.3961303181694973
$ws.Cells.Item(24, 14).Value = 2.652701845052164
$ws.Cells.Item(25, 2).Value = 0.7266699888265862
$ws.Cells.Item(25, 3).Value = 0.08130938494247175
$ws.Cells.Item(25, 4).Value = 0.07616478366076684
$ws.Cells.Item(25, 5).Value = 0.1269188124458012
$ws.Cells.Item(25, 6).Value = 1.856593735678047
$ws.Cells.Item(25, 8).Value = 0.07973214163530429
$ws.Cells.Item(25, 11).Value = 0.7267107997776634
$ws.Cells.Item(25, 13).Value = 0.3452474711761084
$ws.Cells.Item(25, 14).Value = 2.641304851207082
